$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(4, 1).Value = 58.0
$ws.Cells.Item(4, 2).Value = 91.0
$ws.Cells.Item(4, 3).Value = 33.0
$ws.Cells.Item(4, 4).Value = 8.0
$ws.Cells.Item(4, 5).Value = 102.0

$ws.Cells.Item(5, 1).Value = 110.0
$ws.Cells.Item(5, 2).Value = 30.0
$ws.Cells.Item(5, 3).Value = 123.0
$ws.Cells.Item(5, 4).Value = 129.0
$ws.Cells.Item(5, 5).Value = 4.0

$ws.Cells.Item(6, 1).Value = 106.0
$ws.Cells.Item(6, 2).Value = 90.0
$ws.Cells.Item(6, 3).Value = 38.0
$ws.Cells.Item(6, 4).Value = 8.0
$ws.Cells.Item(6, 5).Value = 8.0

$ws.Cells.Item(7, 1).Value = 106.0
$ws.Cells.Item(7, 2).Value = 90.0
$ws.Cells.Item(7, 3).Value = 38.0
$ws.Cells.Item(7, 4).Value = 8.0
$ws.Cells.Item(7, 5).Value = 8.0

$ws.Cells.Item(8, 1).Value = 1.0
$ws.Cells.Item(8, 2).Value = 3.0
$ws.Cells.Item(8, 3).Value = 4.0

$ws.Cells.Item(9, 1).Value = 1.0
$ws.Cells.Item(9, 2).Value = 3.0
$ws.Cells.Item(9, 3).Value = 4.0

$ws.Cells.Item(10, 1).Value = 1.0
$ws.Cells.Item(10, 2).Value = 3.0
$ws.Cells.Item(10, 3).Value = 4.0

$ws.Cells.Item(11, 1).Value = 1.0
$ws.Cells.Item(11, 2).Value = 3.0
$ws.Cells.Item(11, 3).Value = 4.0

$ws.Cells.Item(12, 1).Value = 1.0
$ws.Cells.Item(12, 2).Value = 3.0
$ws.Cells.Item(12, 3).Value = 4.0

$ws.Cells.Item(13, 1).Value = 1.0
$ws.Cells.Item(13, 2).Value = 3.0
$ws.Cells.Item(13, 3).Value = 4.0
